# Fruta / hortaliza, semanal
# Insert a new weekly price record (Primera / Segunda grades, week of
# 2021-09-21) ahead of the existing Kiwi price history, pushing the
# remaining rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows right above the current row 73.
$ws.Rows("73:74").Insert()

# New row 73: Primera grade for the week of 2021-09-21 (serial 44460).
$ws.Range("A73").Value = 11
$ws.Range("B73").Value = "Vega Monumental Concepción"
$ws.Range("C73").Value = "Bíobío"
$ws.Range("D73").Value = 44460
$ws.Range("E73").Value = 8
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100101
$ws.Range("H73").Value = "Berries"
$ws.Range("I73").Value = 100101007
$ws.Range("J73").Value = "Kiwi"
$ws.Range("K73").Value = "Hayward"
$ws.Range("L73").Value = "Primera"
$ws.Range("M73").Value = 200
$ws.Range("N73").Value = 13000
$ws.Range("O73").Value = 14000
$ws.Range("P73").Value = 13500
$ws.Range("Q73").Value = "`$/bandeja 18 kilos"
$ws.Range("R73").Value = "Región de O'Higgins"
$ws.Range("S73").Value = 750
$ws.Range("T73").Value = 18

# New row 74: Segunda grade for the same week.
$ws.Range("A74").Value = 11
$ws.Range("B74").Value = "Vega Monumental Concepción"
$ws.Range("C74").Value = "Bíobío"
$ws.Range("D74").Value = 44460
$ws.Range("E74").Value = 8
$ws.Range("F74").Value = "Fruta"
$ws.Range("G74").Value = 100101
$ws.Range("H74").Value = "Berries"
$ws.Range("I74").Value = 100101007
$ws.Range("J74").Value = "Kiwi"
$ws.Range("K74").Value = "Hayward"
$ws.Range("L74").Value = "Segunda"
$ws.Range("M74").Value = 100
$ws.Range("N74").Value = 11000
$ws.Range("O74").Value = 11000
$ws.Range("P74").Value = 11000
$ws.Range("Q74").Value = "`$/bandeja 18 kilos"
$ws.Range("R74").Value = "Región de O'Higgins"
$ws.Range("S74").Value = 611
$ws.Range("T74").Value = 18
